# Fruta / hortaliza, semanal
#
# A new weekly price record (Piña, Macroferia Regional de Talca) is
# inserted as row 184, pushing the existing rows 184:281 down to 185:282.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 184 (shifts everything below it down).
$ws.Rows("184").Insert()

# Populate the new row with the new record's data.
$ws.Range("A184").Value = 5
$ws.Range("B184").Value = "Macroferia Regional de Talca"
$ws.Range("C184").Value = "Maule"
$ws.Range("D184").Value = 44806
$ws.Range("E184").Value = 7
$ws.Range("F184").Value = "Fruta"
$ws.Range("G184").Value = 100108
$ws.Range("H184").Value = "Tropicales y subtropicales"
$ws.Range("I184").Value = 100108005
$ws.Range("J184").Value = "Piña"
$ws.Range("K184").Value = "Caramelo"
$ws.Range("L184").Value = "Segunda"
$ws.Range("M184").Value = 200
$ws.Range("N184").Value = 18000
$ws.Range("O184").Value = 18000
$ws.Range("P184").Value = 18000
$ws.Range("Q184").Value = "`$/caja 14 unidades"
$ws.Range("R184").Value = "Ecuador"
$ws.Range("S184").Value = 1286
$ws.Range("T184").Value = 14
